# Duty.xlsx — add two new duty rows ("SAS" / "ADD") above the existing
# "Afternoon"/"Night" rows, pushing them down from rows 3-4 to rows 5-6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 3:4 (existing rows 3:4 shift down to 5:6).
# Excel's row-insert carries the formatting of the row above down into the
# new rows (date style on A, centered style on B, text style on D/F/H/J).
$ws.Rows("3:4").Insert()

# The new rows only use columns A-D, so drop the carried-over formatting
# in E:J entirely (no stray empty/styled cells left behind).
$ws.Range("E3:J4").Clear()

# Fill in the date/shift cells first (column A/B) ...
$ws.Range("A3").Value2 = 45759
$ws.Range("B3").Value2 = "Morning"
$ws.Range("A4").Value2 = 45759
$ws.Range("B4").Value2 = "Morning"

# ... then column C for both new rows ...
$ws.Range("C3").Value2 = "SAS"
$ws.Range("C4").Value2 = "ADD"

# ... then column D for both new rows (keeps new shared-string insertion
# order: SAS, ADD, 01255, 424545).
$ws.Range("D3").Value2 = "01255"
$ws.Range("D4").Value2 = "424545"

# Match the saved selection/active cell.
$ws.Range("D8").Select() | Out-Null
